# Auto-generated edit script: updates crafting profit figures (columns H-N)
# across multiple profession sheets, per the scheduled runner data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1159.7778
$ws.Range("I19").Value = 1486
$ws.Range("J19").Value = 1066.5714
$ws.Range("K19").Value = 1486
$ws.Range("L19").Value = 1066.5714
$ws.Range("M19").Value = -1311
$ws.Range("N19").Value = -1416.5714
$ws.Range("H28").Value = 5433.7
$ws.Range("I28").Value = 7617.357
$ws.Range("J28").Value = 338.5
$ws.Range("K28").Value = 7617.357
$ws.Range("L28").Value = 338.5
$ws.Range("M28").Value = -7132.357
$ws.Range("N28").Value = -1308.5
$ws.Range("H33").Value = 1539.75
$ws.Range("I33").Value = 238.21428
$ws.Range("J33").Value = 3361.9
$ws.Range("K33").Value = 238.21428
$ws.Range("L33").Value = 3361.9
$ws.Range("M33").Value = -9.214280000000002
$ws.Range("N33").Value = -3819.9
$ws.Range("H98").Value = 637.71875
$ws.Range("I98").Value = 663.5599999999999
$ws.Range("K98").Value = 663.5599999999999
$ws.Range("M98").Value = 834.4400000000001
$ws.Range("H122").Value = 637.71875
$ws.Range("I122").Value = 663.5599999999999
$ws.Range("K122").Value = 1990.68
$ws.Range("M122").Value = 459.3200000000002
$ws.Range("H141").Value = 4612.125
$ws.Range("I141").Value = 4672.364
$ws.Range("K141").Value = 14017.092
$ws.Range("M141").Value = -8837.091999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 631.3
$ws.Range("I97").Value = 546.0714
$ws.Range("J97").Value = 830.1667
$ws.Range("K97").Value = 546.0714
$ws.Range("L97").Value = 830.1667
$ws.Range("M97").Value = -50.07140000000004
$ws.Range("N97").Value = -1822.1667
$ws.Range("H102").Value = 4187.6216
$ws.Range("I102").Value = 2444.3704
$ws.Range("K102").Value = 2444.3704
$ws.Range("M102").Value = -822.3703999999998
$ws.Range("H132").Value = 1504.4889
$ws.Range("I132").Value = 1500.0476
$ws.Range("J132").Value = 1566.6666
$ws.Range("K132").Value = 4500.142800000001
$ws.Range("L132").Value = 4699.9998
$ws.Range("M132").Value = -1970.142800000001
$ws.Range("N132").Value = -9759.9998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2777.318
$ws.Range("I86").Value = 2143
$ws.Range("J86").Value = 3887.375
$ws.Range("K86").Value = 2143
$ws.Range("L86").Value = 3887.375
$ws.Range("M86").Value = -1020
$ws.Range("N86").Value = -6133.375
$ws.Range("H89").Value = 2777.318
$ws.Range("I89").Value = 2143
$ws.Range("J89").Value = 3887.375
$ws.Range("K89").Value = 10715
$ws.Range("L89").Value = 19436.875
$ws.Range("M89").Value = -5099
$ws.Range("N89").Value = -30668.875
$ws.Range("H134").Value = 2638.0715
$ws.Range("I134").Value = 2544
$ws.Range("J134").Value = 2873.25
$ws.Range("K134").Value = 7632
$ws.Range("L134").Value = 8619.75
$ws.Range("M134").Value = -5097
$ws.Range("N134").Value = -13689.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 50780
$ws.Range("J129").Value = 50780
$ws.Range("L129").Value = 50780
$ws.Range("N129").Value = -60780
$ws.Range("H132").Value = 1357.6857
$ws.Range("I132").Value = 1357.6857
$ws.Range("K132").Value = 4073.0571
$ws.Range("M132").Value = -1543.0571
$ws.Range("H134").Value = 1693.9143
$ws.Range("I134").Value = 1687.8529
$ws.Range("K134").Value = 5063.5587
$ws.Range("M134").Value = -2528.5587

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 95.26667
$ws.Range("I12").Value = 95.5
$ws.Range("J12").Value = 95.111115
$ws.Range("K12").Value = 286.5
$ws.Range("L12").Value = 285.333345
$ws.Range("M12").Value = -113.5
$ws.Range("N12").Value = -631.333345
$ws.Range("H80").Value = 1498.25
$ws.Range("J80").Value = 1499
$ws.Range("L80").Value = 4497
$ws.Range("N80").Value = -6369
$ws.Range("H83").Value = 1498.25
$ws.Range("J83").Value = 1499
$ws.Range("L83").Value = 13491
$ws.Range("N83").Value = -22851
$ws.Range("H87").Value = 3000
$ws.Range("I87").Value = 3000
$ws.Range("K87").Value = 9000
$ws.Range("M87").Value = -7752
$ws.Range("H90").Value = 3000
$ws.Range("I90").Value = 3000
$ws.Range("K90").Value = 27000
$ws.Range("M90").Value = -20760
$ws.Range("H103").Value = 421.7143
$ws.Range("I103").Value = 315
$ws.Range("J103").Value = 501.75
$ws.Range("K103").Value = 945
$ws.Range("L103").Value = 1505.25
$ws.Range("M103").Value = -66
$ws.Range("N103").Value = -3263.25
$ws.Range("H107").Value = 1270.9
$ws.Range("I107").Value = 1813.4166
$ws.Range("J107").Value = 457.125
$ws.Range("K107").Value = 5440.2498
$ws.Range("L107").Value = 1371.375
$ws.Range("M107").Value = -3520.2498
$ws.Range("N107").Value = -5211.375
$ws.Range("H133").Value = 5800
$ws.Range("J133").Value = 5272.727
$ws.Range("L133").Value = 15818.181
$ws.Range("N133").Value = -25938.181

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 15000
$ws.Range("J39").Value = 15000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -16064
$ws.Range("H49").Value = 20479.4
$ws.Range("I49").Value = 14000
$ws.Range("J49").Value = 24799
$ws.Range("K49").Value = 14000
$ws.Range("L49").Value = 24799
$ws.Range("M49").Value = -13816
$ws.Range("N49").Value = -25167
$ws.Range("H122").Value = 3245.606
$ws.Range("I122").Value = 2596
$ws.Range("K122").Value = 7788
$ws.Range("M122").Value = -5338
$ws.Range("H126").Value = 3031
$ws.Range("I126").Value = 2830.1428
$ws.Range("J126").Value = 3499.6667
$ws.Range("K126").Value = 8490.428400000001
$ws.Range("L126").Value = 10499.0001
$ws.Range("M126").Value = -6020.428400000001
$ws.Range("N126").Value = -15439.0001
$ws.Range("H132").Value = 2940
$ws.Range("I132").Value = 2940
$ws.Range("K132").Value = 8820
$ws.Range("M132").Value = -6290

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H7").Value = 2054.7
$ws.Range("J7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("N7").Value = -3224
$ws.Range("H22").Value = 2482.6667
$ws.Range("I22").Value = 1792
$ws.Range("J22").Value = 3864
$ws.Range("K22").Value = 1792
$ws.Range("L22").Value = 3864
$ws.Range("M22").Value = -1497
$ws.Range("N22").Value = -4454
$ws.Range("H27").Value = 2482.6667
$ws.Range("I27").Value = 1792
$ws.Range("J27").Value = 3864
$ws.Range("K27").Value = 1792
$ws.Range("L27").Value = 3864
$ws.Range("M27").Value = -1685
$ws.Range("N27").Value = -4078
$ws.Range("H42").Value = 24997
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 24997
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 24997
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -26123
$ws.Range("H49").Value = 24997
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 24997
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 24997
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -25291
$ws.Range("H61").Value = 3084.2075
$ws.Range("I61").Value = 2276.6562
$ws.Range("K61").Value = 2276.6562
$ws.Range("M61").Value = -2074.6562
$ws.Range("H93").Value = 3593.25
$ws.Range("I93").Value = 726.0625
$ws.Range("J93").Value = 7416.1665
$ws.Range("K93").Value = 726.0625
$ws.Range("L93").Value = 7416.1665
$ws.Range("M93").Value = 521.9375
$ws.Range("N93").Value = -9912.166499999999
$ws.Range("H113").Value = 3084.2075
$ws.Range("I113").Value = 2276.6562
$ws.Range("K113").Value = 2276.6562
$ws.Range("M113").Value = -106.6561999999999
$ws.Range("H122").Value = 7702.0312
$ws.Range("I122").Value = 6265.12
$ws.Range("J122").Value = 12833.857
$ws.Range("K122").Value = 18795.36
$ws.Range("L122").Value = 38501.571
$ws.Range("M122").Value = -16345.36
$ws.Range("N122").Value = -43401.571
$ws.Range("H126").Value = 2054.7
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 2156.0222
$ws.Range("J132").Value = 2108.1082
$ws.Range("L132").Value = 6324.3246
$ws.Range("N132").Value = -11384.3246

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 50250000
$ws.Range("I5").Value = 500000
$ws.Range("K5").Value = 500000
$ws.Range("M5").Value = -499888
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H113").Value = 718.26086
$ws.Range("I113").Value = 386.4
$ws.Range("K113").Value = 1159.2
$ws.Range("M113").Value = 1010.8
$ws.Range("H122").Value = 4622.609
$ws.Range("I122").Value = 1804.5
$ws.Range("K122").Value = 5413.5
$ws.Range("M122").Value = -2963.5
$ws.Range("H136").Value = 2145.158
$ws.Range("I136").Value = 1450.2142
$ws.Range("J136").Value = 4091
$ws.Range("K136").Value = 4350.642599999999
$ws.Range("L136").Value = 12273
$ws.Range("M136").Value = -1800.642599999999
$ws.Range("N136").Value = -17373
